$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 35000
$ws.Range("J21").Value = 35000
$ws.Range("L21").Value = 35000
$ws.Range("N21").Value = -35936
$ws.Range("H23").Value = 35000
$ws.Range("J23").Value = 35000
$ws.Range("L23").Value = 35000
$ws.Range("N23").Value = -35468
$ws.Range("H33").Value = 138
$ws.Range("I33").Value = 138.43243
$ws.Range("K33").Value = 138.43243
$ws.Range("M33").Value = 90.56756999999999
$ws.Range("H62").Value = 4040.7273
$ws.Range("I62").Value = 3700.5557
$ws.Range("J62").Value = 4276.231
$ws.Range("K62").Value = 3700.5557
$ws.Range("L62").Value = 4276.231
$ws.Range("M62").Value = -3076.5557
$ws.Range("N62").Value = -5524.231
$ws.Range("H65").Value = 4040.7273
$ws.Range("I65").Value = 3700.5557
$ws.Range("J65").Value = 4276.231
$ws.Range("K65").Value = 18502.7785
$ws.Range("L65").Value = 21381.155
$ws.Range("M65").Value = -15382.7785
$ws.Range("N65").Value = -27621.155
$ws.Range("H70").Value = 1137
$ws.Range("I70").Value = 1069.75
$ws.Range("J70").Value = 1226.6666
$ws.Range("K70").Value = 3209.25
$ws.Range("L70").Value = 3679.9998
$ws.Range("M70").Value = -2939.25
$ws.Range("N70").Value = -4219.9998
$ws.Range("H73").Value = 1137
$ws.Range("I73").Value = 1069.75
$ws.Range("J73").Value = 1226.6666
$ws.Range("K73").Value = 3209.25
$ws.Range("L73").Value = 3679.9998
$ws.Range("M73").Value = -2273.25
$ws.Range("N73").Value = -5551.9998
$ws.Range("H86").Value = 9832.666999999999
$ws.Range("J86").Value = 18082
$ws.Range("L86").Value = 18082
$ws.Range("N86").Value = -20328
$ws.Range("H89").Value = 9832.666999999999
$ws.Range("J89").Value = 18082
$ws.Range("L89").Value = 90410
$ws.Range("N89").Value = -101642
$ws.Range("H96").Value = 20834972
$ws.Range("I96").Value = 83334240
$ws.Range("J96").Value = 1882.5555
$ws.Range("K96").Value = 250002720
$ws.Range("L96").Value = 5647.666499999999
$ws.Range("M96").Value = -250001347
$ws.Range("N96").Value = -8393.666499999999
$ws.Range("H100").Value = 3212.5
$ws.Range("I100").Value = 2000
$ws.Range("K100").Value = 2000
$ws.Range("M100").Value = -1459
$ws.Range("H112").Value = 3969310.8
$ws.Range("J112").Value = 4274581
$ws.Range("L112").Value = 12823743
$ws.Range("N112").Value = -12825959
$ws.Range("H115").Value = 393.1111
$ws.Range("I115").Value = 393.1111
$ws.Range("K115").Value = 1179.3333
$ws.Range("M115").Value = 387.6667
$ws.Range("H138").Value = 120405.01
$ws.Range("I138").Value = 2281.7273
$ws.Range("J138").Value = 137729.77
$ws.Range("K138").Value = 6845.1819
$ws.Range("L138").Value = 413189.3099999999
$ws.Range("M138").Value = -1705.1819
$ws.Range("N138").Value = -423469.3099999999
$ws.Range("H141").Value = 1884.9642
$ws.Range("I141").Value = 1611.4783
$ws.Range("J141").Value = 3143
$ws.Range("K141").Value = 4834.4349
$ws.Range("L141").Value = 9429
$ws.Range("M141").Value = 345.5650999999998
$ws.Range("N141").Value = -19789

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I74").Value = 26316348
$ws.Range("K74").Value = 26316348
$ws.Range("M74").Value = -26315474
$ws.Range("I77").Value = 26316348
$ws.Range("K77").Value = 131581740
$ws.Range("M77").Value = -131577372

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 41667336
$ws.Range("I64").Value = 83334290
$ws.Range("J64").Value = 384.08334
$ws.Range("K64").Value = 83334290
$ws.Range("L64").Value = 384.08334
$ws.Range("M64").Value = -83334065
$ws.Range("N64").Value = -834.08334
$ws.Range("H67").Value = 41667336
$ws.Range("I67").Value = 83334290
$ws.Range("J67").Value = 384.08334
$ws.Range("K67").Value = 83334290
$ws.Range("L67").Value = 384.08334
$ws.Range("M67").Value = -83333510
$ws.Range("N67").Value = -1944.08334
$ws.Range("H80").Value = 878.3570999999999
$ws.Range("I80").Value = 1263.2727
$ws.Range("J80").Value = 629.2941
$ws.Range("K80").Value = 1263.2727
$ws.Range("L80").Value = 629.2941
$ws.Range("M80").Value = -265.2727
$ws.Range("N80").Value = -2625.2941
$ws.Range("H83").Value = 878.3570999999999
$ws.Range("I83").Value = 1263.2727
$ws.Range("J83").Value = 629.2941
$ws.Range("K83").Value = 6316.363499999999
$ws.Range("L83").Value = 3146.4705
$ws.Range("M83").Value = -1324.363499999999
$ws.Range("N83").Value = -13130.4705
$ws.Range("H99").Value = 1359.5
$ws.Range("I99").Value = 1283.3334
$ws.Range("J99").Value = 1457.4286
$ws.Range("K99").Value = 1283.3334
$ws.Range("L99").Value = 1457.4286
$ws.Range("N99").Value = -4453.4286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3008.589
$ws.Range("I31").Value = 1287.6809
$ws.Range("J31").Value = 6119.4614
$ws.Range("K31").Value = 1287.6809
$ws.Range("L31").Value = 6119.4614
$ws.Range("M31").Value = -992.6809000000001
$ws.Range("N31").Value = -6709.4614
$ws.Range("H34").Value = 3008.589
$ws.Range("I34").Value = 1287.6809
$ws.Range("J34").Value = 6119.4614
$ws.Range("K34").Value = 1287.6809
$ws.Range("L34").Value = 6119.4614
$ws.Range("M34").Value = -1085.6809
$ws.Range("N34").Value = -6523.4614
$ws.Range("H52").Value = 21336.125
$ws.Range("I52").Value = 7709
$ws.Range("J52").Value = 23282.857
$ws.Range("K52").Value = 7709
$ws.Range("L52").Value = 23282.857
$ws.Range("M52").Value = -7415
$ws.Range("N52").Value = -23870.857
$ws.Range("H107").Value = 1734
$ws.Range("I107").Value = 767
$ws.Range("K107").Value = 767
$ws.Range("M107").Value = 1153
$ws.Range("H134").Value = 948.5161000000001
$ws.Range("I134").Value = 779.4138
$ws.Range("J134").Value = 3400.5
$ws.Range("K134").Value = 2338.2414
$ws.Range("L134").Value = 10201.5
$ws.Range("M134").Value = 196.7586000000001
$ws.Range("N134").Value = -15271.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1144.3784
$ws.Range("I5").Value = 813.8929000000001
$ws.Range("K5").Value = 2441.6787
$ws.Range("M5").Value = -2329.6787
$ws.Range("H131").Value = 652.47
$ws.Range("J131").Value = 770.06757
$ws.Range("L131").Value = 2310.20271
$ws.Range("N131").Value = -12390.20271
$ws.Range("H135").Value = 1144.3784
$ws.Range("I135").Value = 813.8929000000001
$ws.Range("K135").Value = 7325.0361
$ws.Range("M135").Value = -4790.0361

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 30000
$ws.Range("J63").Value = 30000
$ws.Range("L63").Value = 30000
$ws.Range("N63").Value = -31372
$ws.Range("H66").Value = 30000
$ws.Range("J66").Value = 30000
$ws.Range("L66").Value = 90000
$ws.Range("N66").Value = -96864
$ws.Range("H102").Value = 3309.6667
$ws.Range("I102").Value = 2830.7856
$ws.Range("J102").Value = 10014
$ws.Range("K102").Value = 2830.7856
$ws.Range("L102").Value = 10014
$ws.Range("M102").Value = -1208.7856
$ws.Range("N102").Value = -13258

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5155.7144
$ws.Range("I40").Value = 4909.737
$ws.Range("J40").Value = 7492.5
$ws.Range("K40").Value = 4909.737
$ws.Range("L40").Value = 7492.5
$ws.Range("M40").Value = -4773.737
$ws.Range("N40").Value = -7764.5
$ws.Range("H100").Value = 1813.8636
$ws.Range("I100").Value = 1273.4546
$ws.Range("J100").Value = 2354.2727
$ws.Range("K100").Value = 1273.4546
$ws.Range("L100").Value = 2354.2727
$ws.Range("M100").Value = -732.4546
$ws.Range("N100").Value = -3436.2727

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 11491.154
$ws.Range("J64").Value = 17586.25
$ws.Range("L64").Value = 17586.25
$ws.Range("N64").Value = -18082.25
$ws.Range("H67").Value = 11491.154
$ws.Range("J67").Value = 17586.25
$ws.Range("L67").Value = 17586.25
$ws.Range("N67").Value = -19302.25
$ws.Range("H113").Value = 1110.2106
$ws.Range("I113").Value = 1383.4166
$ws.Range("J113").Value = 641.8570999999999
$ws.Range("K113").Value = 4150.2498
$ws.Range("L113").Value = 1925.5713
$ws.Range("M113").Value = -1980.2498
$ws.Range("N113").Value = -6265.5713
$ws.Range("H132").Value = 575.3200000000001
$ws.Range("I132").Value = 453.18918
$ws.Range("J132").Value = 922.9231
$ws.Range("K132").Value = 1359.56754
$ws.Range("L132").Value = 2768.7693
$ws.Range("M132").Value = 1170.43246
$ws.Range("N132").Value = -7828.7693
